# Auto-generated edit script applying the Excalibur_Profits.xlsx diff
# Updates FFXIV leve-profit calculation sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR)
# with refreshed market-board price snapshots. All target sheets are plain
# Excel Tables (Table_<SheetName>) over range A1:N141 with cached numeric values
# (no formulas), so each changed cell is written directly via Range.Value.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 13216.368
$ws.Range("I69").Value = 5561.75
$ws.Range("J69").Value = 15257.6
$ws.Range("K69").Value = 16685.25
$ws.Range("L69").Value = 45772.8
$ws.Range("M69").Value = -15811.25
$ws.Range("N69").Value = -47520.8
$ws.Range("H72").Value = 13216.368
$ws.Range("I72").Value = 5561.75
$ws.Range("J72").Value = 15257.6
$ws.Range("K72").Value = 50055.75
$ws.Range("L72").Value = 137318.4
$ws.Range("M72").Value = -45687.75
$ws.Range("N72").Value = -146054.4
$ws.Range("H74").Value = 5598.2
$ws.Range("I74").Value = 4194.778
$ws.Range("K74").Value = 4194.778
$ws.Range("M74").Value = -3258.778
$ws.Range("H77").Value = 5598.2
$ws.Range("I77").Value = 4194.778
$ws.Range("K77").Value = 20973.89
$ws.Range("M77").Value = -16293.89
$ws.Range("H127").Value = 1522.1
$ws.Range("I127").Value = 1522.1
$ws.Range("K127").Value = 4566.299999999999
$ws.Range("M127").Value = 393.7000000000007
$ws.Range("H129").Value = 3381.0908
$ws.Range("I129").Value = 2774.375
$ws.Range("J129").Value = 4999
$ws.Range("K129").Value = 8323.125
$ws.Range("L129").Value = 14997
$ws.Range("M129").Value = -3323.125
$ws.Range("N129").Value = -24997

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5148.5386
$ws.Range("J2").Value = 8055.5
$ws.Range("L2").Value = 8055.5
$ws.Range("N2").Value = -8281.5
$ws.Range("H116").Value = 5148.5386
$ws.Range("J116").Value = 8055.5
$ws.Range("L116").Value = 8055.5
$ws.Range("N116").Value = -12643.5
$ws.Range("H132").Value = 2448.8333
$ws.Range("I132").Value = 1884.7675
$ws.Range("K132").Value = 5654.3025
$ws.Range("M132").Value = -3124.3025

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5148.5386
$ws.Range("J3").Value = 8055.5
$ws.Range("L3").Value = 8055.5
$ws.Range("N3").Value = -8283.5
$ws.Range("H20").Value = 1663.3334
$ws.Range("I20").Value = 1625.4117
$ws.Range("K20").Value = 1625.4117
$ws.Range("M20").Value = -1378.4117
$ws.Range("H105").Value = 3720
$ws.Range("I105").Value = 3354.5715
$ws.Range("K105").Value = 3354.5715
$ws.Range("M105").Value = -1607.5715

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1253925
$ws.Range("J4").Value = 4485.7144
$ws.Range("L4").Value = 4485.7144
$ws.Range("N4").Value = -4709.7144
$ws.Range("H16").Value = 1997.8572
$ws.Range("H22").Value = 773.1579
$ws.Range("I22").Value = 788.3333
$ws.Range("K22").Value = 788.3333
$ws.Range("M22").Value = -438.3333
$ws.Range("H107").Value = 578.4545000000001
$ws.Range("I107").Value = 489.2
$ws.Range("J107").Value = 715.7692
$ws.Range("K107").Value = 489.2
$ws.Range("L107").Value = 715.7692
$ws.Range("M107").Value = 1430.8
$ws.Range("N107").Value = -4555.7692
$ws.Range("H113").Value = 1997.8572

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 1866.6666
$ws.Range("J12").Value = 1475
$ws.Range("L12").Value = 1475
$ws.Range("N12").Value = -1755
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H20").Value = 1714
$ws.Range("I20").Value = 1614
$ws.Range("K20").Value = 1614
$ws.Range("M20").Value = -1369
$ws.Range("H22").Value = 986
$ws.Range("J22").Value = 781.3333
$ws.Range("L22").Value = 781.3333
$ws.Range("N22").Value = -1839.3333
$ws.Range("H70").Value = 16375
$ws.Range("I70").Value = 20331.834
$ws.Range("J70").Value = 4504.5
$ws.Range("K70").Value = 20331.834
$ws.Range("L70").Value = 4504.5
$ws.Range("M70").Value = -20061.834
$ws.Range("N70").Value = -5044.5
$ws.Range("H73").Value = 16375
$ws.Range("I73").Value = 20331.834
$ws.Range("J73").Value = 4504.5
$ws.Range("K73").Value = 20331.834
$ws.Range("L73").Value = 4504.5
$ws.Range("M73").Value = -19395.834
$ws.Range("N73").Value = -6376.5
$ws.Range("H80").Value = 48575.867
$ws.Range("I80").Value = 68110.625
$ws.Range("K80").Value = 68110.625
$ws.Range("M80").Value = -67112.625
$ws.Range("H83").Value = 48575.867
$ws.Range("I83").Value = 68110.625
$ws.Range("K83").Value = 340553.125
$ws.Range("M83").Value = -335561.125
$ws.Range("H97").Value = 2755.4827
$ws.Range("I97").Value = 1670.5294
$ws.Range("J97").Value = 4292.5
$ws.Range("K97").Value = 1670.5294
$ws.Range("L97").Value = 4292.5
$ws.Range("M97").Value = -1174.5294
$ws.Range("N97").Value = -5284.5
$ws.Range("H126").Value = 3926.3157
$ws.Range("I126").Value = 3412.5
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 10237.5
$ws.Range("L126").Value = 12900
$ws.Range("M126").Value = -7767.5
$ws.Range("N126").Value = -17840

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4000
$ws.Range("J2").Value = 4000
$ws.Range("L2").Value = 4000
$ws.Range("N2").Value = -4224
$ws.Range("H7").Value = 3599.7693
$ws.Range("I7").Value = 3399.7
$ws.Range("K7").Value = 3399.7
$ws.Range("M7").Value = -3287.7
$ws.Range("H20").Value = 269700.28
$ws.Range("J20").Value = 323214.28
$ws.Range("L20").Value = 323214.28
$ws.Range("N20").Value = -323666.28
$ws.Range("H22").Value = 1111.875
$ws.Range("J22").Value = 1200
$ws.Range("L22").Value = 1200
$ws.Range("N22").Value = -1790
$ws.Range("H27").Value = 1111.875
$ws.Range("J27").Value = 1200
$ws.Range("L27").Value = 1200
$ws.Range("N27").Value = -1414
$ws.Range("H61").Value = 2688.9312
$ws.Range("I61").Value = 2585.76
$ws.Range("K61").Value = 2585.76
$ws.Range("M61").Value = -2383.76
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240
$ws.Range("H74").Value = 62576.7
$ws.Range("I74").Value = 45979.8
$ws.Range("K74").Value = 45979.8
$ws.Range("M74").Value = -44981.8
$ws.Range("H77").Value = 62576.7
$ws.Range("I77").Value = 45979.8
$ws.Range("K77").Value = 137939.4
$ws.Range("M77").Value = -132947.4
$ws.Range("H82").Value = 1825.75
$ws.Range("I82").Value = 1712
$ws.Range("J82").Value = 2394.5
$ws.Range("K82").Value = 1712
$ws.Range("L82").Value = 2394.5
$ws.Range("M82").Value = -1351
$ws.Range("N82").Value = -3116.5
$ws.Range("H85").Value = 1825.75
$ws.Range("I85").Value = 1712
$ws.Range("J85").Value = 2394.5
$ws.Range("K85").Value = 1712
$ws.Range("L85").Value = 2394.5
$ws.Range("M85").Value = -464
$ws.Range("N85").Value = -4890.5
$ws.Range("H113").Value = 2688.9312
$ws.Range("I113").Value = 2585.76
$ws.Range("K113").Value = 2585.76
$ws.Range("M113").Value = -415.7600000000002
$ws.Range("H126").Value = 3599.7693
$ws.Range("I126").Value = 3399.7
$ws.Range("K126").Value = 10199.1
$ws.Range("M126").Value = -7729.099999999999
$ws.Range("H132").Value = 2878.3438
$ws.Range("I132").Value = 2903.5667
$ws.Range("K132").Value = 8710.7001
$ws.Range("M132").Value = -6180.7001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 213947.3
$ws.Range("I5").Value = 667526
$ws.Range("J5").Value = 19556.428
$ws.Range("K5").Value = 667526
$ws.Range("L5").Value = 19556.428
$ws.Range("M5").Value = -667414
$ws.Range("N5").Value = -19780.428
$ws.Range("H7").Value = 2800
$ws.Range("I7").Value = 2800
$ws.Range("K7").Value = 2800
$ws.Range("M7").Value = -2687
$ws.Range("H45").Value = 7225.5
$ws.Range("J45").Value = 7225.5
$ws.Range("L45").Value = 7225.5
$ws.Range("N45").Value = -8207.5
$ws.Range("H70").Value = 56719.168
$ws.Range("H73").Value = 56719.168
$ws.Range("H75").Value = 59000
$ws.Range("J75").Value = 59000
$ws.Range("L75").Value = 59000
$ws.Range("N75").Value = -60872
$ws.Range("H78").Value = 59000
$ws.Range("J78").Value = 59000
$ws.Range("L78").Value = 177000
$ws.Range("N78").Value = -186360
$ws.Range("H107").Value = 7332.1665
$ws.Range("J107").Value = 5997.5
$ws.Range("L107").Value = 17992.5
$ws.Range("N107").Value = -21832.5
